# Applies the data/view edits captured in the commit diff:
#   - C2 (alpha_distance_range, Max): 9    -> 10
#   - C3 (beta_distance_range,  Max): 8.3  -> 9
#   - B4 (ratio_threshold_range, Min): 0.7  -> 0.8
#   - active cell selection on Sheet1:      C9 -> C8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 9
$ws.Range("B4").Value = 0.8

$ws.Range("C8").Select()
